$wb = $excel.ActiveWorkbook

$wsTestCases = $wb.Worksheets.Item("TestCases")
$wsKeywords  = $wb.Worksheets.Item("Keywords")
$wsTestData  = $wb.Worksheets.Item("TestData")

# --- Keywords sheet: rename the test-case id, fix the locator id ---
$wsKeywords.Range("B2").Value = "DisneyBookAResortTest"
$wsKeywords.Range("B3").Value = "DisneyBookAResortTest"
$wsKeywords.Range("B4").Value = "DisneyBookAResortTest"
$wsKeywords.Range("E4").Value = "findRates_id"

# --- TestData sheet: rename test-case id, switch browser, add SysDate column ---
$wsTestData.Range("A1").Value = "DisneyBookAResortTest"
$wsTestData.Range("B3").Value = "Chrome"

$wsTestData.Range("D2").Copy()
$wsTestData.Range("E2").PasteSpecial(-4122)
$wsTestData.Range("E2").Value = "SysDate"

$wsTestData.Range("D3").NumberFormat = "mm-dd-yy"
$wsTestData.Range("D3").Copy()
$wsTestData.Range("E3").PasteSpecial(-4122)
$wsTestData.Range("D3").Formula = "=TODAY()+10"
$wsTestData.Range("E3").Formula = "=TODAY()"

# --- TestCases sheet: rename test-case id, wrap the data cell ---
$wsTestCases.Range("B2").Value = "DisneyBookAResortTest"
$wsTestCases.Range("C2").WrapText = $true

# --- Selection / active sheet bookkeeping ---
$wsKeywords.Range("K26").Select()
$wsTestCases.Range("C3").Select()
